# Daily attendance processing - reorders the "Recorded By" (column G) entries
# so that the "System" token is listed first among the comma-separated recorders.
#
# Rules (derived from the authoritative diff):
#   - "X, System"                             -> "System, X"             (for X = dnasr281@gmail.com or backup@backdoor.com)
#   - "system, backup@backdoor.com, System"   -> "System, backup@backdoor.com, system"  (swap first/last, case differs only)
#   - Anything else (already starts with "System", contains "admin@admin.com", single values, etc.)
#     is left untouched.
#
# NOTE: string -eq/-ne comparisons in this host are case-insensitive, so we
# cannot detect "did the value change" by comparing old vs new strings (the
# "system, ... , System" -> "System, ..., system" case only differs by case).
# Instead we track an explicit $changed flag set inside the branch that matched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$startRow = $used.Row
$lastRow = $startRow + $used.Rows.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $val = $cell.Value2

    if ($null -eq $val) { continue }

    $parts = $val -split ", "
    $changed = $false
    $newVal = $val

    if ($parts.Count -eq 3 -and $parts[0] -eq "system" -and $parts[1] -eq "backup@backdoor.com" -and $parts[2] -eq "System") {
        $newVal = "System, backup@backdoor.com, system"
        $changed = $true
    }
    elseif ($parts.Count -eq 2 -and $parts[1] -eq "System" -and ($parts[0] -eq "dnasr281@gmail.com" -or $parts[0] -eq "backup@backdoor.com")) {
        $newVal = "System, " + $parts[0]
        $changed = $true
    }

    if ($changed) {
        $cell.Value = $newVal
    }
}
